$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab06")

# --- Update column header labels (row 2): "2011-20" -> "2012-21" ---
$ws.Range("C2").Value = "Youth literacy rate, population 15-24 years, both sexes (%, most recent measure 2012-21)"
$ws.Range("D2").Value = "Youth literacy rate, population 15-24 years, female (%, most recent measure 2012-21)"
$ws.Range("E2").Value = "Youth literacy rate, population 15-24 years, male (%, most recent measure 2012-21)"
$ws.Range("F2").Value = "Youth literacy rate, population 15-24 years, adjusted gender parity index (GPIA, most recent measure 2012-21)"
$ws.Range("G2").Value = "Adult literacy rate, population 15+ years, both sexes (%, most recent measure 2012-21)"
$ws.Range("H2").Value = "Adult literacy rate, population 15+ years, female (%, most recent measure 2012-21)"
$ws.Range("I2").Value = "Adult literacy rate, population 15+ years, male (%, most recent measure 2012-21)"
$ws.Range("J2").Value = "Adult literacy rate, population 15+ years, adjusted gender parity index (GPIA, most recent measure 2012-21)"

# --- Update statistical aggregate values (rows changed due to recalculated regional figures) ---
# Row 63
$ws.Range("C63").Value = 97.189179418604596
$ws.Range("D63").Value = 97.062205813953497
$ws.Range("E63").Value = 97.345552441860505
$ws.Range("F63").Value = 0.99536697674418995
$ws.Range("G63").Value = 92.662818390804603
$ws.Range("H63").Value = 91.252841839080503
$ws.Range("I63").Value = 94.088169770115002
$ws.Range("J63").Value = 0.96371827586207004

# Row 65
$ws.Range("C65").Value = 94.159400000000005
$ws.Range("D65").Value = 93.293322222222201
$ws.Range("E65").Value = 95.0471537037037
$ws.Range("F65").Value = 0.97658814814814998
$ws.Range("G65").Value = 86.209238518518504
$ws.Range("H65").Value = 82.7912888888889
$ws.Range("I65").Value = 89.690868518518499
$ws.Range("J65").Value = 0.90776037037037005

# Row 66
$ws.Range("C66").Value = 90.352451544117699
$ws.Range("D66").Value = 89.338696397058797
$ws.Range("E66").Value = 91.484950588235293
$ws.Range("F66").Value = 0.96642294117646999
$ws.Range("G66").Value = 83.429200802919695
$ws.Range("H66").Value = 80.276602335766498
$ws.Range("I66").Value = 86.746979999999994
$ws.Range("J66").Value = 0.9031500729927

# Row 78
$ws.Range("C78").Value = 99.502804444444493
$ws.Range("D78").Value = 99.581906666666697
$ws.Range("E78").Value = 99.448632222222301
$ws.Range("F78").Value = 1.0013377777777801
$ws.Range("G78").Value = 98.206751111111203
$ws.Range("H78").Value = 98.091163333333398
$ws.Range("I78").Value = 98.341725555555598
$ws.Range("J78").Value = 0.99752666666667

# Row 79
$ws.Range("C79").Value = 99.436281666666702
$ws.Range("D79").Value = 99.507579166666702
$ws.Range("E79").Value = 99.382738333333407
$ws.Range("F79").Value = 1.0012558333333299
$ws.Range("G79").Value = 97.700460833333395
$ws.Range("H79").Value = 97.273721666666702
$ws.Range("I79").Value = 98.146425833333296
$ws.Range("J79").Value = 0.99112333333333003

# Row 83
$ws.Range("C83").Value = 97.038487500000002
$ws.Range("D83").Value = 96.830372941176506
$ws.Range("E83").Value = 97.262228823529497
$ws.Range("F83").Value = 0.99335323529412001
$ws.Range("G83").Value = 92.191270000000003
$ws.Range("H83").Value = 90.748897681159406
$ws.Range("I83").Value = 93.683205217391304
$ws.Range("J83").Value = 0.96146565217391

# Row 85
$ws.Range("C85").Value = 55.929400000000001
$ws.Range("D85").Value = 41.593719999999998
$ws.Range("E85").Value = 71.243939999999995
$ws.Range("F85").Value = 0.58382000000000001
$ws.Range("G85").Value = 37.266039999999997
$ws.Range("H85").Value = 22.60078
$ws.Range("I85").Value = 52.063429999999997
$ws.Range("J85").Value = 0.43409999999999999

# Row 89
$ws.Range("C89").Value = 98.890513333333402
$ws.Range("D89").Value = 98.983066666666701
$ws.Range("E89").Value = 98.804335897436005
$ws.Range("F89").Value = 1.0017630769230801
$ws.Range("G89").Value = 96.078867692307696
$ws.Range("H89").Value = 95.428970256410295
$ws.Range("I89").Value = 96.743538461538506
$ws.Range("J89").Value = 0.98606794871795

# Row 90
$ws.Range("C90").Value = 99.081181999999998
$ws.Range("D90").Value = 99.337048999999993
$ws.Range("E90").Value = 98.910116000000102
$ws.Range("F90").Value = 1.0043325000000001
$ws.Range("G90").Value = 97.457574285714301
$ws.Range("H90").Value = 97.114433809523803
$ws.Range("I90").Value = 97.7116923809524
$ws.Range("J90").Value = 0.99395047619048005

# Row 92
$ws.Range("C92").Value = 86.944335555555597
$ws.Range("D92").Value = 85.237738888888899
$ws.Range("E92").Value = 88.784855555555595
$ws.Range("F92").Value = 0.95121111111111001
$ws.Range("G92").Value = 70.073353333333401
$ws.Range("H92").Value = 63.846510000000002
$ws.Range("I92").Value = 76.613816666666693
$ws.Range("J92").Value = 0.81668444444443999

# Row 94
$ws.Range("C94").Value = 96.545028666666695
$ws.Range("D94").Value = 96.897270666666699
$ws.Range("E94").Value = 96.208965333333396
$ws.Range("F94").Value = 1.0071333333333301
$ws.Range("G94").Value = 91.649180625
$ws.Range("H94").Value = 91.259870000000006
$ws.Range("I94").Value = 92.042691875000003
$ws.Range("J94").Value = 0.98807562500000001

# Row 96
$ws.Range("C96").Value = 95.565889374999998
$ws.Range("D96").Value = 94.462427500000004
$ws.Range("E96").Value = 96.73821375
$ws.Range("F96").Value = 0.96947749999999999
$ws.Range("G96").Value = 89.922199375000005
$ws.Range("H96").Value = 87.222699375000005
$ws.Range("I96").Value = 92.731848749999997
$ws.Range("J96").Value = 0.92005812499999995

# Row 98
$ws.Range("C98").Value = 89.907127857142896
$ws.Range("D98").Value = 88.405255714285701
$ws.Range("E98").Value = 91.454587857142897
$ws.Range("F98").Value = 0.95826500000000003
$ws.Range("G98").Value = 79.009045
$ws.Range("H98").Value = 74.708944285714296
$ws.Range("I98").Value = 83.412341428571395
$ws.Range("J98").Value = 0.87758857142857005

# Row 99
$ws.Range("C99").Value = 67.765094000000005
$ws.Range("D99").Value = 63.235726
$ws.Range("E99").Value = 72.626716999999999
$ws.Range("F99").Value = 0.84274199999999999
$ws.Range("G99").Value = 56.520147000000001
$ws.Range("H99").Value = 48.778283000000002
$ws.Range("I99").Value = 64.598167000000004
$ws.Range("J99").Value = 0.71631199999999995
